# Add a new test-case row (row 37) to the storage_modifiers worksheet,
# mirroring the layout of the existing test-case rows (A=case #,
# B=description, C=given input, D/E=expected/actual output, F=result).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = 3
$ws.Range("B37").Value = "help command"
$ws.Range("C37").Value = "argv[0] -h"
$ws.Range("D37").Value = "enter given inputs"
$ws.Range("E37").Value = "enter given inputs"
$ws.Range("F37").Value = "PASS"

# Match the author's final selection/scroll state after entering the row.
$ws.Range("F38").Select()
